# Updated cryptos list with GitHub Actions
# Refreshes price (column D) and 1h volume change (column E) figures for
# the crypto table, plus swaps the Monero/Arweave rows (50/51).
#
# All of these cells are plain text in the source workbook (t="inlineStr"),
# e.g. "70.346.22" or "  +2.12%  ". Assigning a bare numeric-looking string
# via COM's Range.Value makes Excel auto-coerce it into a real number
# (dropping formatting like the trailing zero in "0.0370"), so purely
# numeric values are entered with a leading apostrophe (forces text) and
# then the cell style is reset to "Normal" so no stray NumberFormat/
# quote-prefix style sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    if ($value -match '^-?\d+(\.\d+)?$') {
        # Looks like a plain number - force text entry, then drop the
        # resulting "quote prefix" style so the cell keeps its original
        # (default) style index.
        $range.Value = "'" + $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range("D2") "70.263.45"
Set-TextValue $ws.Range("E2") "  +1.92%  "

Set-TextValue $ws.Range("D3") "3.950.36"
Set-TextValue $ws.Range("E3") "  +1.93%  "

Set-TextValue $ws.Range("E4") "  +0.18%  "

Set-TextValue $ws.Range("D5") "611.25"
Set-TextValue $ws.Range("E5") "  +1.30%  "

Set-TextValue $ws.Range("D6") "171.14"
Set-TextValue $ws.Range("E6") "  +4.84%  "

Set-TextValue $ws.Range("D7") "3.951.55"
Set-TextValue $ws.Range("E7") "  +2.02%  "

Set-TextValue $ws.Range("E8") "  +0.09%  "

Set-TextValue $ws.Range("E9") "  +1.42%  "

Set-TextValue $ws.Range("E10") "  +1.99%  "

Set-TextValue $ws.Range("E11") "  +2.61%  "

Set-TextValue $ws.Range("E12") "  +2.73%  "

Set-TextValue $ws.Range("E13") "  +5.97%  "

Set-TextValue $ws.Range("D14") "38.57"
Set-TextValue $ws.Range("E14") "  +4.35%  "

Set-TextValue $ws.Range("D15") "4.616.05"
Set-TextValue $ws.Range("E15") "  +2.07%  "

Set-TextValue $ws.Range("D16") "3.930.43"
Set-TextValue $ws.Range("E16") "  +1.41%  "

Set-TextValue $ws.Range("D17") "70.243.59"
Set-TextValue $ws.Range("E17") "  +1.64%  "

Set-TextValue $ws.Range("D18") "7.71"
Set-TextValue $ws.Range("E18") "  +2.04%  "

Set-TextValue $ws.Range("D19") "18.34"
Set-TextValue $ws.Range("E19") "  +6.72%  "

Set-TextValue $ws.Range("E20") "  -1.03%  "

Set-TextValue $ws.Range("D21") "11.17"
Set-TextValue $ws.Range("E21") "  -2.24%  "

Set-TextValue $ws.Range("D22") "500.30"
Set-TextValue $ws.Range("E22") "  +2.87%  "

Set-TextValue $ws.Range("D23") "0.749"
Set-TextValue $ws.Range("E23") "  +3.75%  "

Set-TextValue $ws.Range("D24") "0.0000168"
Set-TextValue $ws.Range("E24") "  +5.50%  "

Set-TextValue $ws.Range("D25") "86.27"
Set-TextValue $ws.Range("E25") "  +2.57%  "

Set-TextValue $ws.Range("E26") "  +3.14%  "

Set-TextValue $ws.Range("D27") "12.45"
Set-TextValue $ws.Range("E27") "  +2.59%  "

Set-TextValue $ws.Range("D28") "10.33"
Set-TextValue $ws.Range("E28") "  +2.51%  "

Set-TextValue $ws.Range("E29") "  +0.11%  "

Set-TextValue $ws.Range("E30") "  +1.24%  "

Set-TextValue $ws.Range("D31") "4.103.47"
Set-TextValue $ws.Range("E31") "  +1.96%  "

Set-TextValue $ws.Range("D32") "2.45"
Set-TextValue $ws.Range("E32") "  +2.73%  "

Set-TextValue $ws.Range("D33") "7.90"
Set-TextValue $ws.Range("E33") "  -0.98%  "

Set-TextValue $ws.Range("D34") "32.46"
Set-TextValue $ws.Range("E34") "  +0.17%  "

Set-TextValue $ws.Range("D35") "3.914.86"
Set-TextValue $ws.Range("E35") "  +2.32%  "

Set-TextValue $ws.Range("E36") "  +2.10%  "

Set-TextValue $ws.Range("D37") "6.20"
Set-TextValue $ws.Range("E37") "  +4.86%  "

Set-TextValue $ws.Range("E38") "  +1.26%  "

Set-TextValue $ws.Range("D39") "0.142"
Set-TextValue $ws.Range("E39") "  +0.90%  "

Set-TextValue $ws.Range("E40") "  +9.61%  "

Set-TextValue $ws.Range("E41") "  +3.73%  "

Set-TextValue $ws.Range("E42") "  +0.25%  "

Set-TextValue $ws.Range("D43") "2.12"
Set-TextValue $ws.Range("E43") "  +6.72%  "

Set-TextValue $ws.Range("D44") "441.26"
Set-TextValue $ws.Range("E44") "  -0.49%  "

Set-TextValue $ws.Range("D45") "48.33"
Set-TextValue $ws.Range("E45") "  -0.43%  "

Set-TextValue $ws.Range("D46") "8.70"
Set-TextValue $ws.Range("E46") "  +3.42%  "

Set-TextValue $ws.Range("E47") "  +0.03%  "

Set-TextValue $ws.Range("D48") "0.000279"
Set-TextValue $ws.Range("E48") "  +23.65%  "

Set-TextValue $ws.Range("D49") "0.0370"
Set-TextValue $ws.Range("E49") "  +3.60%  "

# Row 50/51 swap: Monero moves up to rank 48 (row 50), Arweave drops to
# rank 49 (row 51), with fresh price/volume figures for each.
Set-TextValue $ws.Range("B50") "Monero"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D50") "143.99"
Set-TextValue $ws.Range("E50") "  +0.61%  "

Set-TextValue $ws.Range("B51") "Arweave"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D51") "40.53"
Set-TextValue $ws.Range("E51") "  +4.72%  "
